$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each pair of rows below holds a duplicate stock line (same item name / unit
# rate in column D) whose Barcode/Batch (B), Sale Rate (E), Quantity (F) and
# Value (G) were recorded against the wrong one of the two rows. Swap those
# four columns between the two rows of each pair so the figures line up with
# the correct barcode/batch.
$rowPairs = @(
    @(49,50), @(76,77), @(84,85), @(86,87), @(98,99), @(102,103), @(109,110),
    @(145,146), @(147,148), @(150,151), @(162,163), @(177,178), @(182,183),
    @(225,226), @(232,233), @(253,254), @(369,370), @(374,375), @(408,410),
    @(413,414), @(417,418), @(427,428), @(438,439), @(487,488), @(502,503),
    @(537,538), @(539,540), @(602,603), @(616,617), @(620,621), @(748,749),
    @(750,751), @(780,781), @(782,783), @(805,806), @(807,808), @(831,832),
    @(835,836), @(839,840), @(841,842), @(843,844), @(845,846), @(861,862),
    @(872,873), @(884,885), @(896,897), @(902,903), @(904,905), @(939,940),
    @(946,947), @(977,978)
)

$cols = @("B", "E", "F", "G")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
